# Append two new timesheet entries (rows 34 and 35) to the "Effort" sheet,
# mirroring the existing data/style pattern used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the date-column formatting (style index carrying the ddd dd/mm/yyyy
# number format) from row 33 into the two new rows, without duplicating the
# number format definition.
$ws.Range("A33").Copy()
$ws.Range("A34:A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 34: 16 Jul 2013, 1h, "Setup creation, prerelease sent to Sudar Muthu"
$ws.Range("A34").Value = 41471
$ws.Range("B34").Value = 1
$ws.Range("D34").Value = "Setup creation, prerelease sent to Sudar Muthu"

# Row 35: 17 Jul 2013, 1.25h, "Revision of manual"
$ws.Range("A35").Value = 41472
$ws.Range("B35").Value = 1.25
$ws.Range("D35").Value = "Revision of manual"

# Match the workbook's recorded selection after the edit.
$null = $ws.Range("B35").Select()
